$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Test run CT 01 (row 2) is now passing and has a new timestamp.
$ws.Range("C2").Value = "Passed"
$ws.Range("H2").Value = "05_05_2020--23_36_02 897"

# CT 03 (row 4) was run this time ("Yes") with a new timestamp.
$ws.Range("B4").Value = "Yes"
$ws.Range("H4").Value = "05_05_2020--23_36_56 694"

# CT 05 (row 6) has a new timestamp.
$ws.Range("H6").Value = "05_05_2020--23_37_26 792"

# Move the active selection to B5, as it was when the file was last saved.
$ws.Range("B5").Select() | Out-Null
